$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Tff3 -> Cxcr4 -> ECs) updated TPM-derived values
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2167346666666667
$ws.Range("H2").Value = 0.650204
$ws.Range("M2").Value = 13.35941066666667
$ws.Range("N2").Value = 40.078232
$ws.Range("O2").Value = 0.4925555025958562
$ws.Range("P2").Value = 0.4925555025958562
$ws.Range("Q2").Value = 2.895447417703111
$ws.Range("R2").Value = 26.059026759328
$ws.Range("S2").Value = 0.4925555025958562
$ws.Range("T2").Value = 0.4925555025958562

# Row 3 (ECs -> Tff3 -> Cxcr4 -> FAPs) updated TPM-derived values
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2167346666666667
$ws.Range("H3").Value = 0.650204
$ws.Range("O3").Value = 0.03774352140193379
$ws.Range("P3").Value = 0.03774352140193379
$ws.Range("Q3").Value = 0.2218722174502222
$ws.Range("R3").Value = 1.996849957052
$ws.Range("S3").Value = 0.03774352140193379
$ws.Range("T3").Value = 0.03774352140193379

# Row 4 (ECs -> Tff3 -> Cxcr4 -> MuSCs) updated TPM-derived values
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2167346666666667
$ws.Range("H4").Value = 0.650204
$ws.Range("M4").Value = 12.73953533333333
$ws.Range("N4").Value = 38.218606
$ws.Range("O4").Value = 0.4697009760022101
$ws.Range("P4").Value = 0.46970097600221
$ws.Range("Q4").Value = 2.761098943958222
$ws.Range("R4").Value = 24.849890495624
$ws.Range("S4").Value = 0.4697009760022101
$ws.Range("T4").Value = 0.46970097600221
